{"js": "// Add a centered, underlined \"Figure: ...\" caption to the empty paragraph\n// that immediately follows the UML Class Diagram picture.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph that directly follows the document's (last) inline\n// picture -- this is the empty, centered caption placeholder paragraph.\nlet captionParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const pictures = paragraphs.items[i].inlinePictures;\n  pictures.load(\"items\");\n  await context.sync();\n  if (pictures.items.length > 0 && i + 1 < paragraphs.items.length) {\n    captionParagraph = paragraphs.items[i + 1];\n  }\n}\n\nif (captionParagraph) {\n  // Insert the caption text as a run inside the (currently empty) paragraph.\n  const captionRange = captionParagraph.insertText(\n    \"Figure: UML Class Diagram of an Online Order Processing System\",\n    Word.InsertLocation.replace\n  );\n  captionRange.font.underline = Word.UnderlineType.single;\n  captionRange.font.size = 14;\n\n  // Also underline the paragraph mark itself, matching Word's behavior of\n  // carrying the last-typed run formatting onto the paragraph's own rPr.\n  captionParagraph.font.underline = Word.UnderlineType.single;\n\n  await context.sync();\n}\n", "ps1": "# Add a centered, underlined \"Figure: ...\" caption to the empty paragraph\n# that immediately follows the UML Class Diagram picture.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that directly follows the document's (last) inline\n# picture -- this is the empty, centered caption placeholder paragraph.\n$targetIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.InlineShapes.Count -gt 0) {\n        $targetIndex = $i + 1\n    }\n}\n\n$i = 0\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($i -eq $targetIndex) {\n        $target = $p\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.InsertAfter(\"Figure: UML Class Diagram of an Online Order Processing System\")\n    # Underline + size apply to the whole paragraph range, covering both the\n    # new run's rPr and the paragraph mark's rPr (pPr/rPr).\n    $target.Range.Font.Size = 14\n    $target.Range.Font.Underline = 1\n}\n"}
